$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

$ws1.Range("A9:K39").Copy()
$ws2.Range("A1:K31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws1.Range("A9:K39").Copy()
$ws2.Range("A1:K31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

Write-Host "done"
